{"js": "// Replace the date line and each three-digit \u00f7 one-digit division problem.\n// All target strings are unique in the document, so a direct search+replace\n// for each old->new pair is safe and unambiguous.\nconst replacements = [\n  [\"2024-04-30 Tuesday\", \"2024-05-01 Wednesday\"],\n  [\"764\u00f76=127, 2\", \"363\u00f77=51, 6\"],\n  [\"211\u00f77=30, 1\", \"537\u00f76=89, 3\"],\n  [\"693\u00f79=77, 0\", \"226\u00f78=28, 2\"],\n  [\"532\u00f78=66, 4\", \"368\u00f76=61, 2\"],\n  [\"435\u00f77=62, 1\", \"118\u00f72=59, 0\"],\n  [\"841\u00f77=120, 1\", \"829\u00f78=103, 5\"],\n  [\"489\u00f77=69, 6\", \"146\u00f72=73, 0\"],\n  [\"320\u00f78=40, 0\", \"868\u00f73=289, 1\"],\n  [\"676\u00f78=84, 4\", \"145\u00f72=72, 1\"],\n  [\"800\u00f72=400, 0\", \"138\u00f73=46, 0\"],\n  [\"536\u00f79=59, 5\", \"873\u00f78=109, 1\"],\n  [\"347\u00f78=43, 3\", \"116\u00f75=23, 1\"],\n  [\"147\u00f76=24, 3\", \"245\u00f74=61, 1\"],\n  [\"739\u00f74=184, 3\", \"584\u00f75=116, 4\"],\n  [\"850\u00f77=121, 3\", \"109\u00f78=13, 5\"],\n  [\"255\u00f76=42, 3\", \"147\u00f73=49, 0\"],\n  [\"955\u00f75=191, 0\", \"506\u00f79=56, 2\"],\n  [\"761\u00f72=380, 1\", \"704\u00f78=88, 0\"],\n  [\"800\u00f74=200, 0\", \"874\u00f73=291, 1\"],\n  [\"512\u00f76=85, 2\", \"745\u00f72=372, 1\"],\n  [\"551\u00f79=61, 2\", \"351\u00f76=58, 3\"],\n  [\"976\u00f79=108, 4\", \"607\u00f74=151, 3\"],\n  [\"837\u00f76=139, 3\", \"661\u00f74=165, 1\"],\n  [\"332\u00f77=47, 3\", \"853\u00f75=170, 3\"],\n  [\"450\u00f79=50, 0\", \"594\u00f73=198, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each three-digit / one-digit division problem.\n# Every old string is unique in the document body, so Find/Replace by exact\n# text (match case, no wildcards) maps unambiguously to its replacement.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-30 Tuesday\", \"2024-05-01 Wednesday\"),\n    @(\"764\u00f76=127, 2\", \"363\u00f77=51, 6\"),\n    @(\"211\u00f77=30, 1\", \"537\u00f76=89, 3\"),\n    @(\"693\u00f79=77, 0\", \"226\u00f78=28, 2\"),\n    @(\"532\u00f78=66, 4\", \"368\u00f76=61, 2\"),\n    @(\"435\u00f77=62, 1\", \"118\u00f72=59, 0\"),\n    @(\"841\u00f77=120, 1\", \"829\u00f78=103, 5\"),\n    @(\"489\u00f77=69, 6\", \"146\u00f72=73, 0\"),\n    @(\"320\u00f78=40, 0\", \"868\u00f73=289, 1\"),\n    @(\"676\u00f78=84, 4\", \"145\u00f72=72, 1\"),\n    @(\"800\u00f72=400, 0\", \"138\u00f73=46, 0\"),\n    @(\"536\u00f79=59, 5\", \"873\u00f78=109, 1\"),\n    @(\"347\u00f78=43, 3\", \"116\u00f75=23, 1\"),\n    @(\"147\u00f76=24, 3\", \"245\u00f74=61, 1\"),\n    @(\"739\u00f74=184, 3\", \"584\u00f75=116, 4\"),\n    @(\"850\u00f77=121, 3\", \"109\u00f78=13, 5\"),\n    @(\"255\u00f76=42, 3\", \"147\u00f73=49, 0\"),\n    @(\"955\u00f75=191, 0\", \"506\u00f79=56, 2\"),\n    @(\"761\u00f72=380, 1\", \"704\u00f78=88, 0\"),\n    @(\"800\u00f74=200, 0\", \"874\u00f73=291, 1\"),\n    @(\"512\u00f76=85, 2\", \"745\u00f72=372, 1\"),\n    @(\"551\u00f79=61, 2\", \"351\u00f76=58, 3\"),\n    @(\"976\u00f79=108, 4\", \"607\u00f74=151, 3\"),\n    @(\"837\u00f76=139, 3\", \"661\u00f74=165, 1\"),\n    @(\"332\u00f77=47, 3\", \"853\u00f75=170, 3\"),\n    @(\"450\u00f79=50, 0\", \"594\u00f73=198, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
